$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, centered, bordered) onto the two
# new header cells before writing their text, so I1/J1 share H1's style
# (same as the rest of row 1) instead of minting a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-12 for the new columns I (I0) and J (IF)
$values = @(
    @(1, 3),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(5, 5),
    @(8, 8),
    @(5, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
